# Updated non-tracing input data
# Replace the static N totals (B3, C3) with formulas that sum the
# compartment rows (S, E, I_asym, I_sym, I_sev, R, D) found in rows 16-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Formula = "=B16+B17+B18+B19+B20+B21+B22"
$ws.Range("C3").Formula = "=C16+C17+C18+C19+C20+C21+C22"

# Move the active selection to F16, matching the edited file's last
# cursor position.
$ws.Range("F16").Select()
